$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '38.735.24'
$ws.Range('D3').Value = '2.098.81'
$ws.Range('E3').Value = '  +0.22%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '''227.42'
$ws.Range('E5').Value = '  -0.60%  '
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('D7').Value = '''62.13'
$ws.Range('E7').Value = '  +1.41%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').Value = '''0.390'
$ws.Range('E9').Value = '  +1.88%  '
$ws.Range('E10').Value = '  -0.33%  '
$ws.Range('E11').Value = '  -1.43%  '
$ws.Range('D12').Value = '''15.67'
$ws.Range('E12').Value = '  +5.33%  '
$ws.Range('D13').Value = '2.409.29'
$ws.Range('E13').Value = '  +0.32%  '
$ws.Range('E14').Value = '  -1.49%  '
$ws.Range('D15').Value = '''0.811'
$ws.Range('E15').Value = '  +3.74%  '
$ws.Range('E16').Value = '  +1.11%  '
$ws.Range('D17').Value = '2.137.66'
$ws.Range('E17').Value = '  +1.86%  '
$ws.Range('D18').Value = '38.703.32'
$ws.Range('E18').Value = '  +0.29%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').Value = '''6.13'
$ws.Range('E19').Value = '  +0.54%  '
$ws.Range('B20').Value = 'Litecoin'
$ws.Range('C20').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D20').Value = '''71.63'
$ws.Range('E20').Value = '  +0.84%  '
$ws.Range('E21').Value = '  +0.42%  '
$ws.Range('D22').Value = '''227.71'
$ws.Range('E22').Value = '  +0.55%  '
$ws.Range('D24').Value = '''2.35'
$ws.Range('E24').Value = '  -2.02%  '
$ws.Range('D25').Value = '''2.32'
$ws.Range('E25').Value = '  -0.29%  '
$ws.Range('D26').Value = '''9.63'
$ws.Range('E26').Value = '  +1.83%  '
$ws.Range('D27').Value = '''171.55'
$ws.Range('E27').Value = '  +0.69%  '
$ws.Range('E28').Value = '  +1.81%  '
$ws.Range('E29').Value = '  +3.70%  '
$ws.Range('D30').Value = '''19.30'
$ws.Range('E30').Value = '  +0.79%  '
$ws.Range('E31').Value = '  +7.91%  '
$ws.Range('E32').Value = '  -0.01%  '
$ws.Range('E33').Value = '  +0.97%  '
$ws.Range('D34').Value = '''4.74'
$ws.Range('E34').Value = '  -0.59%  '
$ws.Range('E35').Value = '  +7.08%  '
$ws.Range('D36').Value = '''0.0617'
$ws.Range('E36').Value = '  +1.65%  '
$ws.Range('D37').Value = '''2.39'
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('E38').Value = '  -0.52%  '
$ws.Range('E39').Value = '  -0.02%  '
$ws.Range('E40').Value = '  -2.32%  '
$ws.Range('D41').Value = '''102.61'
$ws.Range('E41').Value = '  +2.40%  '
$ws.Range('E42').Value = '  +2.86%  '
$ws.Range('D43').Value = '1.527.91'
$ws.Range('E43').Value = '  -1.12%  '
$ws.Range('E44').Value = '  +6.58%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = '''7.88'
$ws.Range('E45').Value = '  +2.57%  '
$ws.Range('B46').Value = 'HuobiToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D46').Value = '''2.81'
$ws.Range('E46').Value = '  -0.97%  '
$ws.Range('D47').Value = '''0.0910'
$ws.Range('E47').Value = '  -0.97%  '
$ws.Range('E48').Value = '  -0.55%  '
$ws.Range('E49').Value = '  +1.80%  '
$ws.Range('E50').Value = '  -1.00%  '
$ws.Range('D51').Value = '2.295.60'
$ws.Range('E51').Value = '  +0.25%  '
